# Auto-generated edit script: refresh static price/profit values in Hades_Profits sheets
# (scheduled runner data update - values only, no formulas)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 201
$ws.Range("I18").Value = 201
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 201
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 83
$ws.Range("N18").ClearContents()
$ws.Range("H29").Value = 781.6667
$ws.Range("I29").Value = 138
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 414
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = -133
$ws.Range("N29").Value = -12562
$ws.Range("H38").Value = 535.7692
$ws.Range("I38").Value = 87.72727
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 263.18181
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = 108.81819
$ws.Range("N38").Value = -9744
$ws.Range("H58").Value = 1504.4615
$ws.Range("I58").Value = 222.57143
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 667.71429
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -517.71429
$ws.Range("N58").Value = -9300
$ws.Range("H113").Value = 4900.2856
$ws.Range("J113").Value = 4922.1113
$ws.Range("L113").Value = 4922.1113
$ws.Range("N113").Value = -11430.1113

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4392.75
$ws.Range("I28").Value = 4392.75
$ws.Range("K28").Value = 4392.75
$ws.Range("M28").Value = -4200.75
$ws.Range("H32").Value = 30232.652
$ws.Range("I32").Value = 29711.574
$ws.Range("J32").Value = 31167.53
$ws.Range("K32").Value = 29711.574
$ws.Range("L32").Value = 31167.53
$ws.Range("M32").Value = -29424.574
$ws.Range("N32").Value = -31741.53
$ws.Range("H99").Value = 4392.75
$ws.Range("I99").Value = 4392.75
$ws.Range("K99").Value = 4392.75
$ws.Range("M99").Value = -1397.75
$ws.Range("H122").Value = 5292981
$ws.Range("I122").Value = 1920.6842
$ws.Range("K122").Value = 5762.0526
$ws.Range("M122").Value = -3312.0526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8336074.5
$ws.Range("I134").Value = 2498.9285
$ws.Range("K134").Value = 7496.7855
$ws.Range("M134").Value = -4961.7855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 80000000
$ws.Range("I6").Value = 80000000
$ws.Range("K6").Value = 80000000
$ws.Range("M6").Value = -79999887
$ws.Range("H7").Value = 302.6
$ws.Range("I7").Value = 98
$ws.Range("J7").Value = 353.75
$ws.Range("K7").Value = 98
$ws.Range("L7").Value = 353.75
$ws.Range("M7").Value = 15
$ws.Range("N7").Value = -579.75
$ws.Range("H16").Value = 1804.7693
$ws.Range("I16").Value = 1852.2
$ws.Range("J16").Value = 1775.125
$ws.Range("K16").Value = 1852.2
$ws.Range("L16").Value = 1775.125
$ws.Range("M16").Value = -1565.2
$ws.Range("N16").Value = -2349.125
$ws.Range("H31").Value = 59154.742
$ws.Range("I31").Value = 58399.39
$ws.Range("J31").Value = 59802.19
$ws.Range("K31").Value = 58399.39
$ws.Range("L31").Value = 59802.19
$ws.Range("M31").Value = -58104.39
$ws.Range("N31").Value = -60392.19
$ws.Range("H34").Value = 59154.742
$ws.Range("I34").Value = 58399.39
$ws.Range("J34").Value = 59802.19
$ws.Range("K34").Value = 58399.39
$ws.Range("L34").Value = 59802.19
$ws.Range("M34").Value = -58197.39
$ws.Range("N34").Value = -60206.19
$ws.Range("H38").Value = 900
$ws.Range("I38").Value = 900
$ws.Range("K38").Value = 900
$ws.Range("M38").Value = -523
$ws.Range("H42").Value = 7056
$ws.Range("I42").Value = 7056
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 7056
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("M42").Value = -6463
$ws.Range("H46").Value = 900
$ws.Range("I46").Value = 900
$ws.Range("K46").Value = 900
$ws.Range("M46").Value = -689
$ws.Range("H58").Value = 18183904
$ws.Range("I58").Value = 32260020
$ws.Range("J58").Value = 2255.4167
$ws.Range("K58").Value = 32260020
$ws.Range("L58").Value = 2255.4167
$ws.Range("M58").Value = -32259817
$ws.Range("N58").Value = -2661.4167
$ws.Range("H99").Value = 1666.6666
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -4746
$ws.Range("H113").Value = 1804.7693
$ws.Range("I113").Value = 1852.2
$ws.Range("J113").Value = 1775.125
$ws.Range("K113").Value = 1852.2
$ws.Range("L113").Value = 1775.125
$ws.Range("M113").Value = 317.8
$ws.Range("N113").Value = -6115.125
$ws.Range("H126").Value = 1666.6666
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -10190
$ws.Range("H136").Value = 18183904
$ws.Range("I136").Value = 32260020
$ws.Range("J136").Value = 2255.4167
$ws.Range("K136").Value = 96780060
$ws.Range("L136").Value = 6766.250100000001
$ws.Range("M136").Value = -96777510
$ws.Range("N136").Value = -11866.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 161.2
$ws.Range("I2").Value = 151.83333
$ws.Range("J2").Value = 175.25
$ws.Range("K2").Value = 910.9999799999999
$ws.Range("L2").Value = 1051.5
$ws.Range("M2").Value = -797.9999799999999
$ws.Range("N2").Value = -1277.5
$ws.Range("H4").Value = 3228275.8
$ws.Range("I4").Value = 1821.5714
$ws.Range("J4").Value = 4169325
$ws.Range("K4").Value = 5464.7142
$ws.Range("L4").Value = 12507975
$ws.Range("M4").Value = -5352.7142
$ws.Range("N4").Value = -12508199
$ws.Range("H6").Value = 370.33334
$ws.Range("I6").Value = 64.40000000000001
$ws.Range("J6").Value = 1900
$ws.Range("K6").Value = 193.2
$ws.Range("L6").Value = 5700
$ws.Range("M6").Value = -80.20000000000002
$ws.Range("N6").Value = -5926
$ws.Range("H7").Value = 316.3158
$ws.Range("I7").Value = 194.61539
$ws.Range("J7").Value = 580
$ws.Range("K7").Value = 583.84617
$ws.Range("L7").Value = 1740
$ws.Range("M7").Value = -471.84617
$ws.Range("N7").Value = -1964
$ws.Range("H10").Value = 2540.7222
$ws.Range("I10").Value = 1002.8461
$ws.Range("J10").Value = 6539.2
$ws.Range("K10").Value = 3008.5383
$ws.Range("L10").Value = 19617.6
$ws.Range("M10").Value = -2869.5383
$ws.Range("N10").Value = -19895.6
$ws.Range("H75").Value = 2929.2856
$ws.Range("J75").Value = 3581
$ws.Range("L75").Value = 10743
$ws.Range("N75").Value = -12739
$ws.Range("H78").Value = 2929.2856
$ws.Range("J78").Value = 3581
$ws.Range("L78").Value = 32229
$ws.Range("N78").Value = -42213
$ws.Range("H107").Value = 681.7659
$ws.Range("I107").Value = 635.1
$ws.Range("J107").Value = 764.1177
$ws.Range("K107").Value = 1905.3
$ws.Range("L107").Value = 2292.3531
$ws.Range("M107").Value = 14.69999999999982
$ws.Range("N107").Value = -6132.3531
$ws.Range("H131").Value = 1006.5333
$ws.Range("J131").Value = 1037.5385
$ws.Range("L131").Value = 3112.6155
$ws.Range("N131").Value = -13192.6155
$ws.Range("H132").Value = 1937.2667
$ws.Range("I132").Value = 1116.05
$ws.Range("J132").Value = 2594.24
$ws.Range("K132").Value = 10044.45
$ws.Range("L132").Value = 23348.16
$ws.Range("M132").Value = -7514.449999999999
$ws.Range("N132").Value = -28408.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2192.1365
$ws.Range("I113").Value = 1449.5834
$ws.Range("J113").Value = 3083.2
$ws.Range("K113").Value = 1449.5834
$ws.Range("L113").Value = 3083.2
$ws.Range("M113").Value = 720.4166
$ws.Range("N113").Value = -7423.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 21204.51
$ws.Range("I132").Value = 2208
$ws.Range("J132").Value = 52548.75
$ws.Range("K132").Value = 6624
$ws.Range("L132").Value = 157646.25
$ws.Range("M132").Value = -4094
$ws.Range("N132").Value = -162706.25
